# Add a new header cell D1 with label "CL.RET.CODE:1" to match the
# Corporate Customer excel file layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "CL.RET.CODE:1"

# Move the active selection to the newly added cell, as in the source file.
$ws.Range("D1").Select()
